$wb = $excel.ActiveWorkbook

# --- Schedule sheet updates ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E3").Value = -221.8333455
$wsSchedule.Range("F3").Value = -6.520674470899471
$wsSchedule.Range("B4").Value = 46044.125
$wsSchedule.Range("C4").Value = 6
$wsSchedule.Range("D4").Value = 22.68
$wsSchedule.Range("E4").Value = 694.9715955
$wsSchedule.Range("F4").Value = 30.64248657407408
$wsSchedule.Range("A5").Value = 46044.29166666666
$wsSchedule.Range("C5").Value = 9
$wsSchedule.Range("D5").Value = 34.02
$wsSchedule.Range("E5").Value = -92.05708200000002
$wsSchedule.Range("F5").Value = -2.705969488536156

# --- Detailed sheet updates ---
$wsDetailed = $wb.Worksheets.Item("Detailed")
$wsDetailed.Range("B33").Value = -22.86107
$wsDetailed.Range("B34").Value = -6.8
$wsDetailed.Range("B35").Value = -12.03602
$wsDetailed.Range("C35").Value = "historical"
$wsDetailed.Range("B36").Value = 0.51002
$wsDetailed.Range("C36").Value = "historical"
$wsDetailed.Range("B37").Value = 36.2
$wsDetailed.Range("C37").Value = "historical"
$wsDetailed.Range("B38").Value = 56.98
$wsDetailed.Range("B40").Value = 72.70141
$wsDetailed.Range("B41").Value = 73.2
$wsDetailed.Range("B42").Value = 78
$wsDetailed.Range("B44").Value = 66.37147
$wsDetailed.Range("B45").Value = 61.49545
$wsDetailed.Range("B46").Value = 57.09
$wsDetailed.Range("B47").Value = 57.9708
$wsDetailed.Range("B48").Value = 60.72545
$wsDetailed.Range("B49").Value = 60.93768
$wsDetailed.Range("B51").Value = 62.40053
$wsDetailed.Range("E55").Value = "ON"
$wsDetailed.Range("B56").Value = 57.31
$wsDetailed.Range("B57").Value = 63.73519
$wsDetailed.Range("B58").Value = 65.29086
$wsDetailed.Range("B59").Value = 64.02934
$wsDetailed.Range("B60").Value = 63.77464
$wsDetailed.Range("B61").Value = 73.2
$wsDetailed.Range("B62").Value = 64.29331
$wsDetailed.Range("B63").Value = 57.06
$wsDetailed.Range("E63").Value = "OFF"
$wsDetailed.Range("B64").Value = 26.46866
$wsDetailed.Range("B66").Value = -4.94207
$wsDetailed.Range("B67").Value = -7.71244
$wsDetailed.Range("B68").Value = -13.55778
$wsDetailed.Range("B69").Value = -14
$wsDetailed.Range("B70").Value = -12.01
$wsDetailed.Range("B71").Value = -13.29403
$wsDetailed.Range("B72").Value = -13.05128
$wsDetailed.Range("B73").Value = -7.94077
$wsDetailed.Range("B74").Value = -7.87131
$wsDetailed.Range("B75").Value = -7.81855
$wsDetailed.Range("B76").Value = -7.6547
$wsDetailed.Range("B77").Value = -5.85289
$wsDetailed.Range("B78").Value = -5.81525
$wsDetailed.Range("B79").Value = -0.89511
$wsDetailed.Range("B80").Value = 0.51
$wsDetailed.Range("B82").Value = 0.50922
$wsDetailed.Range("B83").Value = -3.46657
$wsDetailed.Range("B84").Value = -10.86473
$wsDetailed.Range("B85").Value = -7.58517
$wsDetailed.Range("B86").Value = -5.04895
$wsDetailed.Range("B87").Value = -0.41742
$wsDetailed.Range("B88").Value = 10.40351
$wsDetailed.Range("B89").Value = 53.90469
$wsDetailed.Range("B90").Value = 31.99831
$wsDetailed.Range("B91").Value = 51.666
$wsDetailed.Range("B92").Value = 40.99584
$wsDetailed.Range("B93").Value = 50.38784
$wsDetailed.Range("B94").Value = 30.67259
$wsDetailed.Range("B95").Value = 56.85877
$wsDetailed.Range("B96").Value = 56.98
$wsDetailed.Range("B97").Value = 48.40638

Write-Output "Applied run 149 updates to optimisation_result workbook."
